# DDD refactor: duplicate the AHU_2_1 block of rows (2-11) into a new
# AHU_2_2 block (rows 12-21), with device_instance 833 (was 806) and
# object_instance 60 (was 80). All other columns (metric_name, object_type,
# object_instance-index, note) are carried over unchanged from the source
# rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows 12-21, mirroring rows 2-11 for device "AHU_2_2".
$data = @(
    @("AHU_2_2", "VFDFB", 833, "ANALOG_INPUT",  1, 60, "频率反馈"),
    @("AHU_2_2", "CVFB",  833, "ANALOG_INPUT",  2, 60, "水阀开度反馈"),
    @("AHU_2_2", "RT",    833, "ANALOG_INPUT",  3, 60, "回风温度"),
    @("AHU_2_2", "ST",    833, "ANALOG_INPUT",  4, 60, "送风温度"),
    @("AHU_2_2", "S",     833, "BINARY_INPUT",  5, 60, "运行状态"),
    @("AHU_2_2", "AM",    833, "BINARY_INPUT",  6, 60, "自动状态"),
    @("AHU_2_2", "AL",    833, "BINARY_INPUT",  7, 60, "故障状态"),
    @("AHU_2_2", "C",     833, "BINARY_OUTPUT", 0, 60, "启停控制"),
    @("AHU_2_2", "VFDC",  833, "ANALOG_OUTPUT", 0, 60, "频率设定"),
    @("AHU_2_2", "CVC",   833, "ANALOG_OUTPUT", 1, 60, "水阀开度设定")
)

$startRow = 12
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Match the existing table formatting (style index used by the header &
# all other data rows: centered text, "Songti SC Regular" font).
$newRange = $ws.Range("A12:G21")
$newRange.HorizontalAlignment = -4108
$newRange.VerticalAlignment = -4108
$newRange.Font.Name = "Songti SC Regular"

# Match the workbook's recorded selection after the edit.
$ws.Range("D18").Select()
